# Cap nhat tao examSlot
# The exam-room assignment (ExamRoomID, column B) for every student on the
# sheet is updated to the newly generated room/slot "R0025" (value keeps
# the same trailing-space padding the sheet already used for room codes).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2:B6").Value = "R0025     "

# Leave the selection where the user's cursor ended up after the edit.
$ws.Range("D4").Select()
